$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D" = -0.05769999999999999
    "G" = -1.408928571428572
    "H" = -2.232142857142857
    "I" = -1.667560343469096
    "J" = -1.667560343469096
    "K" = -27.6
    "L" = -2.464285714285714
    "U" = 14.7
    "V" = 0.03428171641791045
    "W" = -0.1704756022235948
    "X" = 0.06915004822062507
    "Y" = -0.2396256504442199
    "Z" = 0.06740962536991101
    "AA" = -0.1124096180349719
    "AB" = 0.06352004328524176
    "AC" = -0.1759296613202136
    "AD" = 64.2
    "AE" = 0.1483792342693489
    "AF" = 64.34837923426934
    "AG" = 49.64837923426934
    "AH" = 0.1304848235214431
    "AI" = 0.3564051890533702
    "AJ" = 0.1037695630064185
    "AK" = 0.2993600508096522
    "AL" = 2.68
    "AM" = 2.68
    "AN" = -3.810767495696563
    "AO" = -7.014925373134329
    "AP" = -2.947016040497972
    "AQ" = -7.014925373134329
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
